# RLS (Recursive Least Squares) solved workbook update.
# - Insert 3 helper columns (E:G) on Sheet1 for the RLS model estimate,
#   squared error and a spacer column; this pushes the old "parameter
#   guess" table (old F:G) out to I:J.
# - Add the solved A/B parameter estimates in J5:J6 (with labels in I5:I6).
# - Add the SSE objective cell F34 = SUM(F2:F32).
# - Record the Solver parameters that were used to obtain the solution as
#   hidden, sheet-scoped defined names (exactly what Excel's Solver add-in
#   persists into the workbook once you click "Solve").
# - Re-anchor the chart on Sheet1 so it still sits to the right of the
#   (now wider) data table.
# - Leave the cursor where the user was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Make room for the new E/F/G helper columns (shifts old F/G -> I/J).
# ---------------------------------------------------------------------
$ws.Columns("E:G").Insert()

# ---------------------------------------------------------------------
# 2. Solved RLS parameter estimates (what Solver converged on).
# ---------------------------------------------------------------------
$ws.Range("I5").Value = "A"
$ws.Range("J5").Value = 0.89998914757305348
$ws.Range("I6").Value = "B"
$ws.Range("J6").Value = 0.50003810710633689

# ---------------------------------------------------------------------
# 3. Model estimate column E = (A*C) + (B*B_col), squared error in F,
#    spacer column G left blank (but formatted like F).
# ---------------------------------------------------------------------
$ws.Range("E2:E32").Formula = "=(`$J`$5*C2)+(`$J`$6*B2)"
$ws.Range("E2:E32").Style = "Normal"

$ws.Range("F2:F32").Formula = "=(E2-D2)^2"
$ws.Range("F2:F32").NumberFormat = "0.0"

$ws.Range("G2:G32").NumberFormat = "0.0"
$ws.Range("G2:G32").ClearContents()

# ---------------------------------------------------------------------
# 4. Objective cell Solver minimised: sum of squared errors.
# ---------------------------------------------------------------------
$ws.Range("F34").Formula = "=SUM(F2:F32)"
$ws.Range("F34").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# 5. Persist the Solver run parameters (hidden, sheet-scoped names) —
#    this is exactly what Excel writes after Data > Solver > Solve.
# ---------------------------------------------------------------------
function Add-HiddenName($name, $refersTo) {
    $n = $ws.Names.Add($name, $refersTo)
    $n.Visible = $false
}

Add-HiddenName "solver_adj" "=Sheet1!`$J`$5:`$J`$6"
Add-HiddenName "solver_cvg" "=0.0001"
Add-HiddenName "solver_drv" "=1"
Add-HiddenName "solver_eng" "=1"
Add-HiddenName "solver_est" "=1"
Add-HiddenName "solver_itr" "=2147483647"
Add-HiddenName "solver_mip" "=2147483647"
Add-HiddenName "solver_mni" "=30"
Add-HiddenName "solver_mrt" "=0.075"
Add-HiddenName "solver_msl" "=2"
Add-HiddenName "solver_neg" "=1"
Add-HiddenName "solver_nod" "=2147483647"
Add-HiddenName "solver_num" "=0"
Add-HiddenName "solver_nwt" "=1"
Add-HiddenName "solver_opt" "=Sheet1!`$F`$34"
Add-HiddenName "solver_pre" "=0.000001"
Add-HiddenName "solver_rbv" "=1"
Add-HiddenName "solver_rlx" "=2"
Add-HiddenName "solver_rsd" "=0"
Add-HiddenName "solver_scl" "=1"
Add-HiddenName "solver_sho" "=2"
Add-HiddenName "solver_ssz" "=100"
Add-HiddenName "solver_tim" "=2147483647"
Add-HiddenName "solver_tol" "=0.01"
Add-HiddenName "solver_typ" "=2"
Add-HiddenName "solver_val" "=0"
Add-HiddenName "solver_ver" "=3"

# ---------------------------------------------------------------------
# 6. Re-anchor the chart: it used to start at column I (idx 8); now that
#    E:G were inserted before it, it should start at column L (idx 11).
# ---------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = $ws.Columns("L").Left

# ---------------------------------------------------------------------
# 7. Leave the selection where the user last left it.
# ---------------------------------------------------------------------
$ws.Range("H13").Select()
